$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
# B2 keeps a numeric-looking string "2" (same representation as original "1").
# Force text storage (avoid numeric auto-conversion) via a leading apostrophe,
# then reset the style back to Normal so no stray number-format style is left.
$c = $ws.Range("B2")
$c.Value = "'2"
$c.Style = "Normal"

$ws.Range("G2").Value = -0.6295071740486586
$ws.Range("H2").Value = -0.684341859014348
$ws.Range("I2").Value = -0.8683718028696193
$ws.Range("J2").Value = -0.8683718028696193
$ws.Range("K2").Value = -12.48
$ws.Range("L2").Value = -0.7785402370555209
$ws.Range("X2").Value = 0.05843076096327017
$ws.Range("AB2").Value = 0.05843076096327017
$ws.Range("AL2").Value = 0.143
$ws.Range("AM2").Value = -0.196
$ws.Range("AO2").Value = -97.34265734265735
$ws.Range("AQ2").Value = 71.02040816326529

# --- Row 3 updates ---
$ws.Range("B3").Value = "Kahoot! AS (OB:KAHOT)"

$ws.Range("G3").Value = -0.5344274809160305
$ws.Range("H3").Value = -0.601526717557252
$ws.Range("I3").Value = -0.7381679389312977
$ws.Range("J3").Value = -0.7381679389312977
$ws.Range("K3").Value = -8.08
$ws.Range("L3").Value = -0.6167938931297711
$ws.Range("X3").Value = 0.05843076096327017
$ws.Range("AB3").Value = 0.05843076096327017
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = -0.338
$ws.Range("AO3").ClearContents()
$ws.Range("AQ3").Value = 28.6094674556213

# --- Row 4: new row ---
$ws.Range("A4").Value = "Norway"
$ws.Range("B4").Value = "Play Magnus AS (OB:PMG)"
$ws.Range("C4").Value = "Entertainment"

$ws.Range("G4").Value = -1.054607508532423
$ws.Range("H4").Value = -1.054607508532423
$ws.Range("I4").Value = -1.450511945392491
$ws.Range("J4").Value = -1.450511945392491
$ws.Range("K4").Value = -4.4
$ws.Range("L4").Value = -1.501706484641638
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("X4").Value = 0.05843076096327017
$ws.Range("AB4").Value = 0.05843076096327017
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AJ4").Value = 0
$ws.Range("AL4").Value = 0.143
$ws.Range("AM4").Value = 0.142
$ws.Range("AN4").Value = -0
$ws.Range("AO4").Value = -29.72027972027972
$ws.Range("AP4").Value = -0
$ws.Range("AQ4").Value = -29.92957746478874
